$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (dates) for rows 2-27 from 2023-10-05 (45204) to 2023-10-08 (45207)
$ws.Range("C2:C27").Value = 45207
